# Insert a new worksheet "2023-10-16" as the first sheet in the workbook,
# containing the standard Email/Mobile No./Skills header row and a single
# data row where only the Skills column is populated ("R, C, P").

$wb = $excel.ActiveWorkbook

# Add a brand new worksheet before the current first sheet so it becomes sheet #1.
$firstSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($firstSheet)
$newSheet.Name = "2023-10-16"

# Header row
$newSheet.Range("A1").Value = "Email"
$newSheet.Range("B1").Value = "Mobile No."
$newSheet.Range("C1").Value = "Skills"

# Data row: Email/Mobile left blank, only Skills populated
$newSheet.Range("C2").Value = "R, C, P"

$wb.Worksheets.Item(1).Select()
